$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row after the existing "langue_sigle" rows (row 39 -> new row 40)
$ws.Rows.Item(40).Insert()
$ws.Cells.Item(40, 1).Value = "langue_sigle"
$ws.Cells.Item(40, 2).Value = "ZZ"
$ws.Cells.Item(40, 3).Value = "xyz"

# Insert a new row after the existing "vide" row (now row 43 -> new row 44)
$ws.Rows.Item(44).Insert()
$ws.Cells.Item(44, 1).Value = "vide"
$ws.Cells.Item(44, 2).Value = "missing"
$ws.Cells.Item(44, 3).Value = "manquant"

# Grow the table to include the two new rows
$lo.Resize($ws.Range("A1:C44"))

# Keep view state close to the committed state
$ws.Range("C39").Select()
